$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H133").Value = 38893.5
$ws.Range("J133").Value = 38893.5
$ws.Range("L133").Value = 38893.5
$ws.Range("N133").Value = -49013.5

$ws.Range("H138").Value = 2340.2415
$ws.Range("J138").Value = 2692.7742
$ws.Range("L138").Value = 8078.3226
$ws.Range("N138").Value = -18358.3226

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 51508.6
$ws.Range("I2").Value = 1491.7059
$ws.Range("K2").Value = 1491.7059
$ws.Range("M2").Value = -1378.7059

$ws.Range("H74").Value = 2158.6
$ws.Range("I74").Value = 1274.8235
$ws.Range("J74").Value = 7166.6665
$ws.Range("K74").Value = 1274.8235
$ws.Range("L74").Value = 7166.6665
$ws.Range("M74").Value = -400.8235
$ws.Range("N74").Value = -8914.666499999999

$ws.Range("H77").Value = 2158.6
$ws.Range("I77").Value = 1274.8235
$ws.Range("J77").Value = 7166.6665
$ws.Range("K77").Value = 6374.1175
$ws.Range("L77").Value = 35833.3325
$ws.Range("M77").Value = -2006.1175
$ws.Range("N77").Value = -44569.3325

$ws.Range("H86").Value = 37475.75
$ws.Range("J86").Value = 46539.332
$ws.Range("L86").Value = 46539.332
$ws.Range("N86").Value = -48911.332

$ws.Range("H89").Value = 37475.75
$ws.Range("J89").Value = 46539.332
$ws.Range("L89").Value = 139617.996
$ws.Range("N89").Value = -151473.996

$ws.Range("H116").Value = 51508.6
$ws.Range("I116").Value = 1491.7059
$ws.Range("K116").Value = 1491.7059
$ws.Range("M116").Value = 802.2941000000001

$ws.Range("H122").Value = 1773.5555
$ws.Range("I122").Value = 1938.3846
$ws.Range("J122").Value = 1345
$ws.Range("K122").Value = 5815.1538
$ws.Range("L122").Value = 4035
$ws.Range("M122").Value = -3365.1538
$ws.Range("N122").Value = -8935

$ws.Range("H132").Value = 2847.4
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 2847.4
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 8542.200000000001
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -13602.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 51508.6
$ws.Range("I3").Value = 1491.7059
$ws.Range("K3").Value = 1491.7059
$ws.Range("M3").Value = -1377.7059

$ws.Range("H99").Value = 1109.3043
$ws.Range("I99").Value = 830.3333
$ws.Range("J99").Value = 1632.375
$ws.Range("K99").Value = 830.3333
$ws.Range("L99").Value = 1632.375
$ws.Range("M99").Value = 667.6667
$ws.Range("N99").Value = -4628.375

$ws.Range("H105").Value = 81865.8
$ws.Range("I105").Value = 112953.22
$ws.Range("J105").Value = 64379.125
$ws.Range("K105").Value = 112953.22
$ws.Range("L105").Value = 64379.125
$ws.Range("M105").Value = -111206.22
$ws.Range("N105").Value = -67873.125

$ws.Range("H107").Value = 76959200
$ws.Range("I107").Value = 125057790
$ws.Range("K107").Value = 125057790
$ws.Range("M107").Value = -125055870

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 1834.4615
$ws.Range("I86").Value = 1375.8823
$ws.Range("J86").Value = 2700.6667
$ws.Range("K86").Value = 1375.8823
$ws.Range("L86").Value = 2700.6667
$ws.Range("M86").Value = -252.8823
$ws.Range("N86").Value = -4946.6667

$ws.Range("H89").Value = 1834.4615
$ws.Range("I89").Value = 1375.8823
$ws.Range("J89").Value = 2700.6667
$ws.Range("K89").Value = 6879.4115
$ws.Range("L89").Value = 13503.3335
$ws.Range("M89").Value = -1263.4115
$ws.Range("N89").Value = -24735.3335

$ws.Range("H103").Value = 4120
$ws.Range("I103").Value = 4120
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 4120
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = -2948
$ws.Range("N103").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 4500
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 4500
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 13500
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = -13724

$ws.Range("H17").Value = 850.1667
$ws.Range("I17").Value = 419.8
$ws.Range("K17").Value = 1259.4
$ws.Range("M17").Value = -1090.4

$ws.Range("H34").Value = 348
$ws.Range("J34").Value = 582
$ws.Range("L34").Value = 1746
$ws.Range("N34").Value = -1914

$ws.Range("H39").Value = 1434.6666
$ws.Range("J39").Value = 1902
$ws.Range("L39").Value = 5706
$ws.Range("N39").Value = -6294

$ws.Range("H131").Value = 845.96
$ws.Range("J131").Value = 854.5417
$ws.Range("L131").Value = 2563.6251
$ws.Range("N131").Value = -12643.6251

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 125002060
$ws.Range("I80").Value = 500002500
$ws.Range("J80").Value = 1917.6666
$ws.Range("K80").Value = 500002500
$ws.Range("L80").Value = 1917.6666
$ws.Range("M80").Value = -500001502
$ws.Range("N80").Value = -3913.6666

$ws.Range("H83").Value = 125002060
$ws.Range("I83").Value = 500002500
$ws.Range("J83").Value = 1917.6666
$ws.Range("K83").Value = 2500012500
$ws.Range("L83").Value = 9588.333000000001
$ws.Range("M83").Value = -2500007508
$ws.Range("N83").Value = -19572.333

$ws.Range("H102").Value = 3278.2856
$ws.Range("I102").Value = 3201.4546
$ws.Range("K102").Value = 3201.4546
$ws.Range("M102").Value = -1579.4546

$ws.Range("H122").Value = 759.5789
$ws.Range("I122").Value = 708.375
$ws.Range("K122").Value = 2125.125
$ws.Range("M122").Value = 324.875

$ws.Range("H132").Value = 2708.5
$ws.Range("I132").Value = 2275.1765
$ws.Range("J132").Value = 3760.8572
$ws.Range("K132").Value = 6825.529500000001
$ws.Range("L132").Value = 11282.5716
$ws.Range("M132").Value = -4295.529500000001
$ws.Range("N132").Value = -16342.5716

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2844.9443
$ws.Range("I7").Value = 1822.1111
$ws.Range("J7").Value = 3867.7778
$ws.Range("K7").Value = 1822.1111
$ws.Range("L7").Value = 3867.7778
$ws.Range("M7").Value = -1710.1111
$ws.Range("N7").Value = -4091.7778

$ws.Range("H40").Value = 2354.3635
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()

$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()

$ws.Range("H122").Value = 3887.375
$ws.Range("I122").Value = 5099.8
$ws.Range("J122").Value = 1866.6666
$ws.Range("K122").Value = 15299.4
$ws.Range("L122").Value = 5599.9998
$ws.Range("M122").Value = -12849.4
$ws.Range("N122").Value = -10499.9998

$ws.Range("H126").Value = 2844.9443
$ws.Range("I126").Value = 1822.1111
$ws.Range("J126").Value = 3867.7778
$ws.Range("K126").Value = 5466.3333
$ws.Range("L126").Value = 11603.3334
$ws.Range("M126").Value = -2996.3333
$ws.Range("N126").Value = -16543.3334
